$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. A new day of data (11/04/2021, serial 44308) has come in. Insert it as a
#    new row 123 by copying row 122 (which is currently the highlighted
#    "latest day" row, using the green "Good" cell style) down one row. This
#    both creates the row and carries its formulas down with references
#    shifted correctly, and gives row 123 the "latest day" styling.
# ---------------------------------------------------------------------------
$ws.Rows(122).Copy()
$ws.Rows(123).Insert()

# ---------------------------------------------------------------------------
# 2. Row 122 is no longer the latest day, so restyle it like a normal data
#    row (matching row 121's yellow "Neutral" styling), keeping its own
#    formulas/values untouched.
# ---------------------------------------------------------------------------
$ws.Range("A121:N121").Copy()
$ws.Range("A122:N122").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Fill in row 123's real data for 11/04/2021; formulas (C/D/E/F/I/K/L)
#    were already carried down correctly by the row insert above and will
#    recalculate automatically once the inputs below are set.
# ---------------------------------------------------------------------------
$ws.Range("A123").Value = 44308
$ws.Range("B123").Value = 3874
$ws.Range("G123").Value = 14306
$ws.Range("H123").Value = 28612
$ws.Range("J123").Value = 7153
$ws.Range("M123").Value = 476352
$ws.Range("N123").Value = 476352

# ---------------------------------------------------------------------------
# 4. Apply the data corrections to doses-administered figures that came in
#    with this upload (11/02 and 11/04 counts were revised).
# ---------------------------------------------------------------------------
$ws.Range("B120").Value = 3831
$ws.Range("B122").Value = 3700

# ---------------------------------------------------------------------------
# 5. The "population left" running total in column D (rows 117-122) is a
#    simple fill-down pattern (D[n] = D[n-1] - B[n]); re-enter it across the
#    whole block so it is stored as a single shared formula, the way it ends
#    up after the new row is filled in.
# ---------------------------------------------------------------------------
$ws.Range("D117:D122").Formula = "=(D116-B117)"

# ---------------------------------------------------------------------------
# 6. Leave the selection where Excel would after this kind of edit.
# ---------------------------------------------------------------------------
$ws.Range("I129").Select()
